# Applies two changes described by the commit diff:
#
# 1. Slide 16's table (the plenary "missing gaps" table) switches its
#    DrawingML table style from the custom "{A423CD0E-...}" style to the
#    built-in "{970F872F-...}" style.
#
# 2. The deck's theme is swapped from the "Integral" color palette back to
#    the default "Office" color palette (the slide master's theme part is
#    rewritten with the stock Office theme colors -- dk1/lt1/dk2/lt2/accent1-6
#    /hlink/folHlink -- in place of the Integral greens/golds).

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 16 -------------------------------------------
$slide = $p.Slides.Item(16)
$tableShape = $slide.Shapes.Item(3)
$table = $tableShape.Table
$table.ApplyStyle("{970F872F-D1AA-43DF-8E26-68B83CF843D6}")

# --- 2. Restore the default Office color theme -----------------------------
# .RGB uses the VBA/OLE COLORREF (BGR) integer encoding: val = R + G*256 +
# B*65536. Decimal values below are that encoding of the stock "Office"
# theme hex colors (dk1=000000, lt1=FFFFFF, dk2=44546A, lt2=E7E6E6,
# accent1=5B9BD5, accent2=ED7D31, accent3=A5A5A5, accent4=FFC000,
# accent5=4472C4, accent6=70AD47, hlink=0563C1, folHlink=954F72).
$colors = $p.SlideMaster.ColorScheme
$colors.Colors(1).RGB  = 0          # dk1      000000
$colors.Colors(2).RGB  = 16777215   # lt1      FFFFFF
$colors.Colors(3).RGB  = 6968388    # dk2      44546A
$colors.Colors(4).RGB  = 15132391   # lt2      E7E6E6
$colors.Colors(5).RGB  = 13998939   # accent1  5B9BD5
$colors.Colors(6).RGB  = 3243501    # accent2  ED7D31
$colors.Colors(7).RGB  = 10855845   # accent3  A5A5A5
$colors.Colors(8).RGB  = 49407      # accent4  FFC000
$colors.Colors(9).RGB  = 12874308   # accent5  4472C4
$colors.Colors(10).RGB = 4697456    # accent6  70AD47
$colors.Colors(11).RGB = 12673797   # hlink    0563C1
$colors.Colors(12).RGB = 7491477    # folHlink 954F72
